# Applies the cryptocurrency price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value.
# "AsText" = $true forces text storage for values that otherwise look
# like plain numbers (Excel would auto-convert "2.430" -> 2.43, etc.),
# matching the inline-string cells already used throughout the sheet.
$edits = @(
    @{ Ref = "D2"; Val = "27.574.94"; AsText = $false },
    @{ Ref = "E2"; Val = "  +1.04%  "; AsText = $false },
    @{ Ref = "D3"; Val = "1.762.79"; AsText = $false },
    @{ Ref = "E3"; Val = "  -1.40%  "; AsText = $false },
    @{ Ref = "D4"; Val = "1.007"; AsText = $true },
    @{ Ref = "E4"; Val = "  +0.39%  "; AsText = $false },
    @{ Ref = "D5"; Val = "336.49"; AsText = $true },
    @{ Ref = "E5"; Val = "  +0.55%  "; AsText = $false },
    @{ Ref = "D6"; Val = "1.004"; AsText = $true },
    @{ Ref = "E6"; Val = "  +0.50%  "; AsText = $false },
    @{ Ref = "D7"; Val = "0.3836"; AsText = $true },
    @{ Ref = "E7"; Val = "  +1.30%  "; AsText = $false },
    @{ Ref = "D8"; Val = "0.3399"; AsText = $true },
    @{ Ref = "E8"; Val = "  -1.32%  "; AsText = $false },
    @{ Ref = "D9"; Val = "47.16"; AsText = $true },
    @{ Ref = "E9"; Val = "  -2.68%  "; AsText = $false },
    @{ Ref = "E10"; Val = "  -5.44%  "; AsText = $false },
    @{ Ref = "D11"; Val = "0.07399"; AsText = $true },
    @{ Ref = "E11"; Val = "  -1.43%  "; AsText = $false },
    @{ Ref = "D12"; Val = "1.006"; AsText = $true },
    @{ Ref = "E12"; Val = "  +0.54%  "; AsText = $false },
    @{ Ref = "D13"; Val = "21.77"; AsText = $true },
    @{ Ref = "E13"; Val = "  -0.75%  "; AsText = $false },
    @{ Ref = "D14"; Val = "6.335"; AsText = $true },
    @{ Ref = "E14"; Val = "  -2.43%  "; AsText = $false },
    @{ Ref = "D15"; Val = "1.763.72"; AsText = $false },
    @{ Ref = "E15"; Val = "  -1.48%  "; AsText = $false },
    @{ Ref = "D16"; Val = "7.026"; AsText = $true },
    @{ Ref = "E16"; Val = "  -1.20%  "; AsText = $false },
    @{ Ref = "E17"; Val = "  -2.64%  "; AsText = $false },
    @{ Ref = "D18"; Val = "0.06666"; AsText = $true },
    @{ Ref = "E18"; Val = "  +0.05%  "; AsText = $false },
    @{ Ref = "D19"; Val = "82.18"; AsText = $true },
    @{ Ref = "E19"; Val = "  -2.23%  "; AsText = $false },
    @{ Ref = "D20"; Val = "1.005"; AsText = $true },
    @{ Ref = "E20"; Val = "  +0.47%  "; AsText = $false },
    @{ Ref = "D21"; Val = "17.29"; AsText = $true },
    @{ Ref = "E21"; Val = "  -0.49%  "; AsText = $false },
    @{ Ref = "D22"; Val = "6.409"; AsText = $true },
    @{ Ref = "E22"; Val = "  -3.61%  "; AsText = $false },
    @{ Ref = "D23"; Val = "27.575.36"; AsText = $false },
    @{ Ref = "E23"; Val = "  +1.02%  "; AsText = $false },
    @{ Ref = "D24"; Val = "12.07"; AsText = $true },
    @{ Ref = "E24"; Val = "  -2.91%  "; AsText = $false },
    @{ Ref = "D25"; Val = "2.379"; AsText = $true },
    @{ Ref = "E25"; Val = "  -1.37%  "; AsText = $false },
    @{ Ref = "D26"; Val = "20.67"; AsText = $true },
    @{ Ref = "E26"; Val = "  -3.41%  "; AsText = $false },
    @{ Ref = "B27"; Val = "ImmutableX"; AsText = $false },
    @{ Ref = "C27"; Val = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; AsText = $false },
    @{ Ref = "D27"; Val = "1.422"; AsText = $true },
    @{ Ref = "E27"; Val = "  -5.03%  "; AsText = $false },
    @{ Ref = "B28"; Val = "LidoDAOToken"; AsText = $false },
    @{ Ref = "C28"; Val = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; AsText = $false },
    @{ Ref = "D28"; Val = "2.430"; AsText = $true },
    @{ Ref = "E28"; Val = "  -5.14%  "; AsText = $false },
    @{ Ref = "D29"; Val = "153.18"; AsText = $true },
    @{ Ref = "E29"; Val = "  -0.37%  "; AsText = $false },
    @{ Ref = "D30"; Val = "134.66"; AsText = $true },
    @{ Ref = "E30"; Val = "  +0.31%  "; AsText = $false },
    @{ Ref = "D31"; Val = "1.965.55"; AsText = $false },
    @{ Ref = "E31"; Val = "  -1.28%  "; AsText = $false },
    @{ Ref = "D32"; Val = "6.088"; AsText = $true },
    @{ Ref = "E32"; Val = "  -0.78%  "; AsText = $false },
    @{ Ref = "D33"; Val = "3.961"; AsText = $true },
    @{ Ref = "E33"; Val = "  -1.46%  "; AsText = $false },
    @{ Ref = "D34"; Val = "0.08727"; AsText = $true },
    @{ Ref = "E34"; Val = "  +0.26%  "; AsText = $false },
    @{ Ref = "D35"; Val = "12.64"; AsText = $true },
    @{ Ref = "E35"; Val = "  -5.11%  "; AsText = $false },
    @{ Ref = "D36"; Val = "0.02412"; AsText = $true },
    @{ Ref = "E36"; Val = "  +2.91%  "; AsText = $false },
    @{ Ref = "D37"; Val = "5.355"; AsText = $true },
    @{ Ref = "E37"; Val = "  -2.26%  "; AsText = $false },
    @{ Ref = "D38"; Val = "0.6736"; AsText = $true },
    @{ Ref = "E38"; Val = "  -3.62%  "; AsText = $false },
    @{ Ref = "B39"; Val = "Algorand"; AsText = $false },
    @{ Ref = "C39"; Val = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; AsText = $false },
    @{ Ref = "D39"; Val = "0.2186"; AsText = $true },
    @{ Ref = "E39"; Val = "  -1.28%  "; AsText = $false },
    @{ Ref = "B40"; Val = "Hedera"; AsText = $false },
    @{ Ref = "C40"; Val = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; AsText = $false },
    @{ Ref = "D40"; Val = "0.06270"; AsText = $true },
    @{ Ref = "E40"; Val = "  -1.35%  "; AsText = $false },
    @{ Ref = "D41"; Val = "1.530"; AsText = $true },
    @{ Ref = "E41"; Val = "  -8.14%  "; AsText = $false },
    @{ Ref = "D42"; Val = "1.247"; AsText = $true },
    @{ Ref = "E42"; Val = "  +0.03%  "; AsText = $false },
    @{ Ref = "D43"; Val = "8.359"; AsText = $true },
    @{ Ref = "E43"; Val = "  -5.67%  "; AsText = $false },
    @{ Ref = "D44"; Val = "14.25"; AsText = $true },
    @{ Ref = "E44"; Val = "  -1.54%  "; AsText = $false },
    @{ Ref = "D45"; Val = "1.004"; AsText = $true },
    @{ Ref = "E45"; Val = "  +0.54%  "; AsText = $false },
    @{ Ref = "D46"; Val = "0.6221"; AsText = $true },
    @{ Ref = "E46"; Val = "  -5.19%  "; AsText = $false },
    @{ Ref = "D47"; Val = "3.828"; AsText = $true },
    @{ Ref = "E47"; Val = "  -0.49%  "; AsText = $false },
    @{ Ref = "D48"; Val = "131.45"; AsText = $true },
    @{ Ref = "E48"; Val = "  +1.54%  "; AsText = $false },
    @{ Ref = "D49"; Val = "2.073"; AsText = $true },
    @{ Ref = "E49"; Val = "  -3.79%  "; AsText = $false },
    @{ Ref = "D50"; Val = "0.07329"; AsText = $true },
    @{ Ref = "E50"; Val = "  +2.41%  "; AsText = $false },
    @{ Ref = "D51"; Val = "1.232"; AsText = $true },
    @{ Ref = "E51"; Val = "  +1.72%  "; AsText = $false }
)

foreach ($e in $edits) {
    $cell = $ws.Range($e.Ref)
    if ($e.AsText) {
        $cell.NumberFormat = "@"
        $cell.Value = $e.Val
        $cell.ClearFormats()
    } else {
        $cell.Value = $e.Val
    }
}
